$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Fill in the missing "STG -> FND" usage for the CMMS JSON pipeline (row 5),
# matching the value already used in row 4 (Dev: SPVB_SAHANA/SPVB_SAHANA_STG2FND)
$ws.Range("G5").Value = "Dev: SPVB_SAHANA/SPVB_SAHANA_STG2FND"

# Row 5 grows taller (like row 4 already is) to fit the wrapped text that now fills G5
$ws.Rows("5").RowHeight = 45

# Update the active selection to the cell that was just completed
$ws.Range("G5").Select()
